# GMS Data Release 1
# Update the data dictionary "panels_applied" sheet:
#   - row 3 (rare_diseases_family_id) field name renamed to "referral_id"
#   - row 6 (sample_id) field name renamed to "platekey"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "referral_id"
$ws.Range("B6").Value = "platekey"
